# Append " (Changed main)" after the existing sentence in the first
# paragraph, as three separate text pieces so they land in their own
# <w:r> runs (matching the target diff) instead of being coalesced into
# the previous run.
#
# Plain Range.InsertAfter() calls get merged into the preceding run when
# the run formatting is identical. Turning on TrackRevisions while
# inserting keeps each insertion in its own run (wrapped in <w:ins>), and
# then accepting each recorded revision individually unwraps the <w:ins>
# markup while leaving the runs split apart - exactly the shape the diff
# expects, without introducing any explicit run formatting.

$d = $word.ActiveDocument

$d.TrackRevisions = $true

$p1 = $d.Paragraphs(1).Range
# Paragraph.Range.End includes the trailing paragraph mark; back up one
# character so the insertion point sits right after "document." and
# before the pilcrow.
$insertionPoint = $p1.End - 1
$r = $d.Range($insertionPoint, $insertionPoint)

$r.InsertAfter(" (")
$r.Collapse(0)
$r.InsertAfter("Changed main")
$r.Collapse(0)
$r.InsertAfter(")")

$d.TrackRevisions = $false

# Accept each recorded revision individually (rather than
# Document.AcceptAllRevisions, which also forces a full repagination that
# strips unrelated <w:lastRenderedPageBreak/> hints elsewhere in the
# document).
for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
    $d.Revisions($i).Accept()
}
